$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-02-20 Thursday"; new = "2025-02-21 Friday"},
    @{old = "33×72=2376"; new = "35×79=2765"},
    @{old = "72×65=4680"; new = "45×98=4410"},
    @{old = "55×86=4730"; new = "43×32=1376"},
    @{old = "45×67=3015"; new = "14×99=1386"},
    @{old = "92×12=1104"; new = "90×52=4680"},
    @{old = "81×31=2511"; new = "92×40=3680"},
    @{old = "79×50=3950"; new = "94×97=9118"},
    @{old = "36×93=3348"; new = "33×19=627"},
    @{old = "41×62=2542"; new = "11×92=1012"},
    @{old = "62×16=992"; new = "30×73=2190"},
    @{old = "37×89=3293"; new = "27×45=1215"},
    @{old = "70×12=840"; new = "64×55=3520"},
    @{old = "84×74=6216"; new = "12×73=876"},
    @{old = "64×92=5888"; new = "69×24=1656"},
    @{old = "18×61=1098"; new = "41×90=3690"},
    @{old = "29×76=2204"; new = "38×33=1254"},
    @{old = "17×44=748"; new = "34×38=1292"},
    @{old = "42×15=630"; new = "84×43=3612"},
    @{old = "75×53=3975"; new = "66×96=6336"},
    @{old = "73×32=2336"; new = "59×69=4071"},
    @{old = "83×20=1660"; new = "98×50=4900"},
    @{old = "37×18=666"; new = "53×40=2120"},
    @{old = "48×70=3360"; new = "36×67=2412"},
    @{old = "58×76=4408"; new = "34×99=3366"},
    @{old = "77×12=924"; new = "53×30=1590"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Host "Done applying $($replacements.Count) replacements"
